$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"22.42000000000007"
$ws.Range("H2").Value = [double]"1.981960426045148e-05"
$ws.Range("I2").Value = [double]"1.981960426045148e-05"
$ws.Range("L2").Value = [double]"59.05951211380641"
$ws.Range("M2").Value = "[34.19850579560922, 83.9205184320036]"
$ws.Range("N2").Value = [double]"1.877359626090147e-05"
$ws.Range("O2").Value = [double]"1.877359626090147e-05"
$ws.Range("P2").Value = [double]"1.742184514603347"
$ws.Range("Q2").Value = "[1.2138686329185768, 2.270500396288118]"
$ws.Range("R2").Value = [double]"3.469432119374005e-08"
$ws.Range("S2").Value = [double]"3.469432119374005e-08"
$ws.Range("T2").Value = [double]"66.84039887302761"
$ws.Range("U2").Value = "[51.25404932770047, 82.42674841835475]"
$ws.Range("V2").Value = [double]"4.15085743554755e-11"
$ws.Range("W2").Value = [double]"4.15085743554755e-11"
$ws.Range("X2").Value = [double]"16.2034434434435"
$ws.Range("Y2").Value = [double]"14.31827827827832"
$ws.Range("Z2").Value = [double]"18.08860860860867"
$ws.Range("F3").Value = [double]"22.42000000000007"
$ws.Range("H3").Value = [double]"0.0002622415560687985"
$ws.Range("I3").Value = [double]"0.0002622415560687985"
$ws.Range("L3").Value = [double]"46.3685727756158"
$ws.Range("M3").Value = "[18.18714639686206, 74.54999915436954]"
$ws.Range("N3").Value = [double]"0.001822701307429231"
$ws.Range("O3").Value = [double]"0.001822701307429231"
$ws.Range("P3").Value = [double]"2.106974051957119"
$ws.Range("Q3").Value = "[1.4906055233248878, 2.7233425805893496]"
$ws.Range("R3").Value = [double]"1.509658065224073e-08"
$ws.Range("S3").Value = [double]"1.509658065224073e-08"
$ws.Range("T3").Value = [double]"57.16187149599351"
$ws.Range("U3").Value = "[42.42348832563974, 71.90025466634728]"
$ws.Range("V3").Value = [double]"6.490814552506663e-10"
$ws.Range("W3").Value = [double]"6.490814552506663e-10"
$ws.Range("X3").Value = [double]"14.90178178178182"
$ws.Range("Y3").Value = [double]"12.70242242242246"
$ws.Range("Z3").Value = [double]"17.10114114114119"
$ws.Range("F4").Value = [double]"22.42000000000007"
$ws.Range("H4").Value = [double]"0.02132371256284349"
$ws.Range("I4").Value = [double]"0.02132371256284349"
$ws.Range("L4").Value = [double]"35.43827454216042"
$ws.Range("M4").Value = "[2.6204908642923783, 68.25605822002846]"
$ws.Range("N4").Value = [double]"0.03493045669497064"
$ws.Range("O4").Value = [double]"0.03493045669497064"
$ws.Range("P4").Value = [double]"2.660447832769735"
$ws.Range("Q4").Value = "[1.3773949772495793, 3.9435006882898906]"
$ws.Range("R4").Value = [double]"0.0001341798557361606"
$ws.Range("S4").Value = [double]"0.0001341798557361606"
$ws.Range("T4").Value = [double]"61.39661508228335"
$ws.Range("U4").Value = "[43.55892258643678, 79.23430757812991]"
$ws.Range("V4").Value = [double]"1.283371497962094e-08"
$ws.Range("W4").Value = [double]"1.283371497962094e-08"
$ws.Range("X4").Value = [double]"12.92684684684688"
$ws.Range("Y4").Value = [double]"8.348588588588612"
$ws.Range("Z4").Value = [double]"17.50510510510516"
$ws.Range("F5").Value = [double]"22.42000000000007"
$ws.Range("H5").Value = [double]"0.007431744497329862"
$ws.Range("I5").Value = [double]"0.007431744497329862"
$ws.Range("L5").Value = [double]"38.93503183556402"
$ws.Range("M5").Value = "[9.459456809239342, 68.4106068618887]"
$ws.Range("N5").Value = [double]"0.01077535047978695"
$ws.Range("O5").Value = [double]"0.01077535047978695"
$ws.Range("P5").Value = [double]"3.025237370123504"
$ws.Range("Q5").Value = "[2.031500354573579, 4.018974385673429]"
$ws.Range("R5").Value = [double]"1.991128184553048e-07"
$ws.Range("S5").Value = [double]"1.991128184553048e-07"
$ws.Range("T5").Value = [double]"68.06344991015791"
$ws.Range("U5").Value = "[50.80702503213385, 85.31987478818198]"
$ws.Range("V5").Value = [double]"4.157676425364798e-10"
$ws.Range("W5").Value = [double]"4.157676425364798e-10"
$ws.Range("X5").Value = [double]"11.62518518518522"
$ws.Range("Y5").Value = [double]"8.079279279279302"
$ws.Range("Z5").Value = [double]"15.17109109109114"
$ws.Range("B6").Value = [double]"0"
$ws.Range("F6").Value = [double]"24.72000000000043"
$ws.Range("H6").Value = [double]"0.120998073104335"
$ws.Range("I6").Value = [double]"0.120998073104335"
$ws.Range("L6").Value = [double]"24.7928930036514"
$ws.Range("M6").Value = "[-4.248146979877085, 53.83393298717988]"
$ws.Range("N6").Value = [double]"0.09239951261996548"
$ws.Range("O6").Value = [double]"0.09239951261996548"
$ws.Range("P6").Value = [double]"2.987500521431735"
$ws.Range("Q6").Value = "[1.3019212798660398, 4.67307976299743]"
$ws.Range("R6").Value = [double]"0.0008636887738302512"
$ws.Range("S6").Value = [double]"0.0008636887738302512"
$ws.Range("T6").Value = [double]"65.02046807232892"
$ws.Range("U6").Value = "[48.45067443988829, 81.59026170476956]"
$ws.Range("V6").Value = [double]"4.766242955867028e-10"
$ws.Range("W6").Value = [double]"4.766242955867028e-10"
$ws.Range("X6").Value = [double]"12.96624624624647"
$ws.Range("Y6").Value = [double]"6.334654654654761"
$ws.Range("Z6").Value = [double]"19.59783783783818"
$ws.Range("F7").Value = [double]"24.72000000000043"
$ws.Range("H7").Value = [double]"0.0001219126133685799"
$ws.Range("I7").Value = [double]"0.0001219126133685799"
$ws.Range("L7").Value = [double]"51.74444881295538"
$ws.Range("M7").Value = "[22.02700634712818, 81.46189127878259]"
$ws.Range("N7").Value = [double]"0.001039871008395199"
$ws.Range("O7").Value = [double]"0.001039871008395199"
$ws.Range("P7").Value = [double]"2.371131992799504"
$ws.Range("Q7").Value = "[1.842816111114733, 2.899447874484274]"
$ws.Range("R7").Value = [double]"1.114530689960702e-11"
$ws.Range("S7").Value = [double]"1.114530689960702e-11"
$ws.Range("T7").Value = [double]"69.89371792657892"
$ws.Range("U7").Value = "[54.17881554974821, 85.60862030340962]"
$ws.Range("V7").Value = [double]"1.45292666786645e-11"
$ws.Range("W7").Value = [double]"1.45292666786645e-11"
$ws.Range("X7").Value = [double]"15.3912312312315"
$ws.Range("Y7").Value = [double]"13.3126726726729"
$ws.Range("Z7").Value = [double]"17.46978978979009"
$ws.Range("F8").Value = [double]"24.72000000000043"
$ws.Range("H8").Value = [double]"1.676083752566448e-06"
$ws.Range("I8").Value = [double]"1.676083752566448e-06"
$ws.Range("L8").Value = [double]"65.03530747565681"
$ws.Range("M8").Value = "[37.49878451232502, 92.5718304389886]"
$ws.Range("N8").Value = [double]"2.057341204575458e-05"
$ws.Range("O8").Value = [double]"2.057341204575458e-05"
$ws.Range("P8").Value = [double]"2.673026782333658"
$ws.Range("Q8").Value = "[2.232763547596349, 3.113290017070966]"
$ws.Range("R8").Value = [double]"6.661338147750939e-16"
$ws.Range("S8").Value = [double]"6.661338147750939e-16"
$ws.Range("T8").Value = [double]"64.27071050745879"
$ws.Range("U8").Value = "[48.962492969189654, 79.57892804572792]"
$ws.Range("V8").Value = [double]"7.546696600968517e-11"
$ws.Range("W8").Value = [double]"7.546696600968517e-11"
$ws.Range("X8").Value = [double]"14.20348348348373"
$ws.Range("Y8").Value = [double]"12.47135135135157"
$ws.Range("Z8").Value = [double]"15.93561561561589"
$ws.Range("F9").Value = [double]"24.72000000000043"
$ws.Range("H9").Value = [double]"0.0002525334600129714"
$ws.Range("I9").Value = [double]"0.0002525334600129714"
$ws.Range("L9").Value = [double]"56.38186576130744"
$ws.Range("M9").Value = "[21.91138138704332, 90.85235013557156]"
$ws.Range("N9").Value = [double]"0.001927623296138936"
$ws.Range("O9").Value = [double]"0.001927623296138936"
$ws.Range("P9").Value = [double]"2.647868883205812"
$ws.Range("Q9").Value = "[2.0440793041375027, 3.251658462274121]"
$ws.Range("R9").Value = [double]"2.186495429157276e-11"
$ws.Range("S9").Value = [double]"2.186495429157276e-11"
$ws.Range("T9").Value = [double]"73.58457628550727"
$ws.Range("U9").Value = "[55.72192920007582, 91.44722337093872]"
$ws.Range("V9").Value = [double]"1.279025774181264e-10"
$ws.Range("W9").Value = [double]"1.279025774181264e-10"
$ws.Range("X9").Value = [double]"14.30246246246271"
$ws.Range("Y9").Value = [double]"11.92696696696717"
$ws.Range("Z9").Value = [double]"16.67795795795825"
$ws.Range("F10").Value = [double]"24.72000000000043"
$ws.Range("H10").Value = [double]"0.0005906754493449906"
$ws.Range("I10").Value = [double]"0.0005906754493449906"
$ws.Range("L10").Value = [double]"45.31484200002466"
$ws.Range("M10").Value = "[16.202258506580343, 74.42742549346897]"
$ws.Range("N10").Value = [double]"0.003023933833587789"
$ws.Range("O10").Value = [double]"0.003023933833587789"
$ws.Range("P10").Value = [double]"2.547237286694427"
$ws.Range("Q10").Value = "[1.9057108589343486, 3.188763714454505]"
$ws.Range("R10").Value = [double]"3.479849741694352e-10"
$ws.Range("S10").Value = [double]"3.479849741694352e-10"
$ws.Range("T10").Value = [double]"56.50105483989638"
$ws.Range("U10").Value = "[41.15589348044735, 71.84621619934542]"
$ws.Range("V10").Value = [double]"2.472387850005475e-09"
$ws.Range("W10").Value = [double]"2.472387850005475e-09"
$ws.Range("X10").Value = [double]"14.69837837837863"
$ws.Range("Y10").Value = [double]"12.17441441441462"
$ws.Range("Z10").Value = [double]"17.22234234234264"
$ws.Range("F11").Value = [double]"24.72000000000043"
$ws.Range("H11").Value = [double]"8.912587465670363e-05"
$ws.Range("I11").Value = [double]"8.912587465670363e-05"
$ws.Range("L11").Value = [double]"51.1160449176632"
$ws.Range("M11").Value = "[21.6736514005847, 80.55843843474169]"
$ws.Range("N11").Value = [double]"0.001071640462158463"
$ws.Range("O11").Value = [double]"0.001071640462158463"
$ws.Range("P11").Value = [double]"2.232763547596349"
$ws.Range("Q11").Value = "[1.7170266154755014, 2.7485004797171966]"
$ws.Range("R11").Value = [double]"3.166955586664244e-11"
$ws.Range("S11").Value = [double]"3.166955586664244e-11"
$ws.Range("T11").Value = [double]"63.55410526357332"
$ws.Range("U11").Value = "[48.282092047469035, 78.8261184796776]"
$ws.Range("V11").Value = [double]"9.657763477832759e-11"
$ws.Range("W11").Value = [double]"9.657763477832759e-11"
$ws.Range("X11").Value = [double]"15.93561561561589"
$ws.Range("Y11").Value = [double]"13.90654654654679"
$ws.Range("Z11").Value = [double]"17.964684684685"
$ws.Range("F12").Value = [double]"24.72000000000043"
$ws.Range("H12").Value = [double]"3.210967635447481e-06"
$ws.Range("I12").Value = [double]"3.210967635447481e-06"
$ws.Range("L12").Value = [double]"57.36621799923068"
$ws.Range("M12").Value = "[33.50605182222415, 81.2263841762372]"
$ws.Range("N12").Value = [double]"1.551532844090708e-05"
$ws.Range("O12").Value = [double]"1.551532844090708e-05"
$ws.Range("P12").Value = [double]"1.754763464167271"
$ws.Range("Q12").Value = "[1.2516054816103468, 2.257921446724196]"
$ws.Range("R12").Value = [double]"9.38234112624059e-09"
$ws.Range("S12").Value = [double]"9.38234112624059e-09"
$ws.Range("T12").Value = [double]"70.92511730913867"
$ws.Range("U12").Value = "[56.68464577873671, 85.16558883954063]"
$ws.Range("V12").Value = [double]"4.716227408607665e-13"
$ws.Range("W12").Value = [double]"4.716227408607665e-13"
$ws.Range("X12").Value = [double]"17.81621621621652"
$ws.Range("Y12").Value = [double]"15.83663663663691"
$ws.Range("Z12").Value = [double]"19.79579579579614"
$ws.Range("F13").Value = [double]"24.72000000000043"
$ws.Range("H13").Value = [double]"5.676803356780002e-05"
$ws.Range("I13").Value = [double]"5.676803356780002e-05"
$ws.Range("L13").Value = [double]"56.67126961539194"
$ws.Range("M13").Value = "[30.925433602000737, 82.41710562878313]"
$ws.Range("N13").Value = [double]"5.903104950721882e-05"
$ws.Range("O13").Value = [double]"5.903104950721882e-05"
$ws.Range("P13").Value = [double]"1.603816069400195"
$ws.Range("Q13").Value = "[0.9874475407679641, 2.220184598032426]"
$ws.Range("R13").Value = [double]"4.109943662733428e-06"
$ws.Range("S13").Value = [double]"4.109943662733428e-06"
$ws.Range("T13").Value = [double]"56.87412632474198"
$ws.Range("U13").Value = "[40.32307863545758, 73.42517401402638]"
$ws.Range("V13").Value = [double]"1.334476817937968e-08"
$ws.Range("W13").Value = [double]"1.334476817937968e-08"
$ws.Range("X13").Value = [double]"18.41009009009041"
$ws.Range("Y13").Value = [double]"15.98510510510538"
$ws.Range("Z13").Value = [double]"20.83507507507543"
$ws.Range("F14").Value = [double]"24.72000000000043"
$ws.Range("H14").Value = [double]"0.0005146943131483228"
$ws.Range("I14").Value = [double]"0.0005146943131483228"
$ws.Range("L14").Value = [double]"53.86818573635328"
$ws.Range("M14").Value = "[23.934832137367877, 83.80153933533867]"
$ws.Range("N14").Value = [double]"0.0007335860877972777"
$ws.Range("O14").Value = [double]"0.0007335860877972777"
$ws.Range("P14").Value = [double]"1.46544762419704"
$ws.Range("Q14").Value = "[0.7484474990534231, 2.182447749340657]"
$ws.Range("R14").Value = [double]"0.0001620006941285279"
$ws.Range("S14").Value = [double]"0.0001620006941285279"
$ws.Range("T14").Value = [double]"66.98397585222456"
$ws.Range("U14").Value = "[48.68063263775656, 85.28731906669256]"
$ws.Range("V14").Value = [double]"2.880695681994894e-09"
$ws.Range("W14").Value = [double]"2.880695681994894e-09"
$ws.Range("X14").Value = [double]"18.9544744744748"
$ws.Range("Y14").Value = [double]"16.13357357357385"
$ws.Range("Z14").Value = [double]"21.77537537537576"
